# daily auto push: 2025-10-10 07:26 UTC
# Append the new daily log entry as the next row (row 89) of the Sheet1
# data table: 2025/10/10, 金, 16, 201.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$newRow = $lastRow + 1

# Column A holds the date as literal text (e.g. "2025/10/10"), matching
# every other row in the sheet -- not a real Excel date. Force the cell to
# Text format before assigning so the date-like string isn't auto-converted
# into a date serial number, then drop the temporary format again so the
# cell is left with the sheet's normal (default) styling, same as its
# neighbours.
$dateCell = $ws.Cells.Item($newRow, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025/10/10"
$dateCell.ClearFormats()

$ws.Cells.Item($newRow, 2).Value = "金"
$ws.Cells.Item($newRow, 3).Value = 16
$ws.Cells.Item($newRow, 4).Value = 201
